$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Merge a new "hand" test-run row into the statistics table ---
# Row 8 previously existed only as a blank, pre-formatted placeholder row
# (it already carries the date number-format + border styling inherited
# from the table). Fill it in with the new run's figures.
$ws.Range("A8").Value = "Apr 16, 2018"
$ws.Range("B8").Value = 21
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 24

# --- Chart: extend the value-axis ceiling so the taller stack still fits ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$valAx = $chart.Axes(2)
$valAx.MaximumScale = 25

# --- Chart: resize/reposition it (it was made taller, starting higher up) ---
$co.Left = 361.4931640625
$co.Top = 68.1
$co.Width = 376.89179687499995
$co.Height = 387.70007874015755

# --- Chart: keep the plot area's inner proportions matching the new box ---
$pa = $chart.PlotArea
$pa.InsideTop = 0.24651669342061266
$pa.InsideHeight = 0.72063222512286429

# --- Selection left where the user's cursor ended up ---
$ws.Range("L11").Select()
